$d = $word.ActiveDocument

$pairs = @(
    @("70×97=6790", "45×86=3870"),
    @("38×61=2318", "89×23=2047"),
    @("73×36=2628", "98×40=3920"),
    @("74×70=5180", "54×40=2160"),
    @("58×51=2958", "20×52=1040"),
    @("11×94=1034", "23×13=299"),
    @("97×27=2619", "71×85=6035"),
    @("38×46=1748", "28×21=588"),
    @("23×85=1955", "78×50=3900"),
    @("62×13=806",  "47×77=3619"),
    @("47×80=3760", "47×34=1598"),
    @("57×81=4617", "43×32=1376"),
    @("52×90=4680", "94×39=3666"),
    @("47×87=4089", "53×98=5194"),
    @("38×14=532",  "38×60=2280"),
    @("74×58=4292", "85×63=5355"),
    @("29×97=2813", "34×32=1088"),
    @("16×28=448",  "65×13=845"),
    @("29×76=2204", "30×44=1320"),
    @("67×36=2412", "33×95=3135"),
    @("18×23=414",  "67×49=3283"),
    @("98×43=4214", "64×25=1600"),
    @("61×38=2318", "76×80=6080"),
    @("95×90=8550", "47×32=1504"),
    @("28×75=2100", "97×16=1552")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
